{"js": "// Replace the date line and each \"a\u00f7b=c, d\" division-table cell with its\n// updated values. Every \"old\" text below is unique within the document, so\n// a plain case-sensitive search+replace for each pair is unambiguous.\nconst replacements = [\n  [\"2026-02-28 Saturday\", \"2026-03-01 Sunday\"],\n  [\"14\u00f73=4, 2\", \"33\u00f75=6, 3\"],\n  [\"85\u00f73=28, 1\", \"10\u00f79=1, 1\"],\n  [\"42\u00f75=8, 2\", \"22\u00f78=2, 6\"],\n  [\"88\u00f79=9, 7\", \"65\u00f76=10, 5\"],\n  [\"41\u00f79=4, 5\", \"30\u00f73=10, 0\"],\n  [\"82\u00f79=9, 1\", \"19\u00f75=3, 4\"],\n  [\"60\u00f79=6, 6\", \"94\u00f74=23, 2\"],\n  [\"95\u00f79=10, 5\", \"62\u00f77=8, 6\"],\n  [\"57\u00f78=7, 1\", \"73\u00f76=12, 1\"],\n  [\"95\u00f77=13, 4\", \"63\u00f74=15, 3\"],\n  [\"21\u00f74=5, 1\", \"72\u00f72=36, 0\"],\n  [\"97\u00f78=12, 1\", \"40\u00f77=5, 5\"],\n  [\"19\u00f77=2, 5\", \"52\u00f73=17, 1\"],\n  [\"26\u00f73=8, 2\", \"27\u00f75=5, 2\"],\n  [\"16\u00f72=8, 0\", \"32\u00f76=5, 2\"],\n  [\"36\u00f75=7, 1\", \"17\u00f79=1, 8\"],\n  [\"29\u00f75=5, 4\", \"66\u00f74=16, 2\"],\n  [\"59\u00f78=7, 3\", \"89\u00f79=9, 8\"],\n  [\"46\u00f76=7, 4\", \"39\u00f79=4, 3\"],\n  [\"76\u00f79=8, 4\", \"43\u00f73=14, 1\"],\n  [\"99\u00f75=19, 4\", \"29\u00f77=4, 1\"],\n  [\"85\u00f75=17, 0\", \"43\u00f73=14, 1\"],\n  [\"36\u00f79=4, 0\", \"80\u00f78=10, 0\"],\n  [\"41\u00f75=8, 1\", \"49\u00f77=7, 0\"],\n  [\"96\u00f73=32, 0\", \"68\u00f72=34, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=c, d\" division-table cell with its\n# updated values. Every \"old\" text below is unique within the document, so\n# a single Find/Replace (wdReplaceOne) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-28 Saturday\", \"2026-03-01 Sunday\"),\n    @(\"14\u00f73=4, 2\", \"33\u00f75=6, 3\"),\n    @(\"85\u00f73=28, 1\", \"10\u00f79=1, 1\"),\n    @(\"42\u00f75=8, 2\", \"22\u00f78=2, 6\"),\n    @(\"88\u00f79=9, 7\", \"65\u00f76=10, 5\"),\n    @(\"41\u00f79=4, 5\", \"30\u00f73=10, 0\"),\n    @(\"82\u00f79=9, 1\", \"19\u00f75=3, 4\"),\n    @(\"60\u00f79=6, 6\", \"94\u00f74=23, 2\"),\n    @(\"95\u00f79=10, 5\", \"62\u00f77=8, 6\"),\n    @(\"57\u00f78=7, 1\", \"73\u00f76=12, 1\"),\n    @(\"95\u00f77=13, 4\", \"63\u00f74=15, 3\"),\n    @(\"21\u00f74=5, 1\", \"72\u00f72=36, 0\"),\n    @(\"97\u00f78=12, 1\", \"40\u00f77=5, 5\"),\n    @(\"19\u00f77=2, 5\", \"52\u00f73=17, 1\"),\n    @(\"26\u00f73=8, 2\", \"27\u00f75=5, 2\"),\n    @(\"16\u00f72=8, 0\", \"32\u00f76=5, 2\"),\n    @(\"36\u00f75=7, 1\", \"17\u00f79=1, 8\"),\n    @(\"29\u00f75=5, 4\", \"66\u00f74=16, 2\"),\n    @(\"59\u00f78=7, 3\", \"89\u00f79=9, 8\"),\n    @(\"46\u00f76=7, 4\", \"39\u00f79=4, 3\"),\n    @(\"76\u00f79=8, 4\", \"43\u00f73=14, 1\"),\n    @(\"99\u00f75=19, 4\", \"29\u00f77=4, 1\"),\n    @(\"85\u00f75=17, 0\", \"43\u00f73=14, 1\"),\n    @(\"36\u00f79=4, 0\", \"80\u00f78=10, 0\"),\n    @(\"41\u00f75=8, 1\", \"49\u00f77=7, 0\"),\n    @(\"96\u00f73=32, 0\", \"68\u00f72=34, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
